$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 150991.8164123378
$ws.Range("J2").Value = 322916.2541075106
$ws.Range("K2").Value = 175139.8326550492
$ws.Range("L2").Value = 344977.530106144

$ws.Range("I3").Value = 619340.4399868009
$ws.Range("J3").Value = 904148.1533773908
$ws.Range("K3").Value = 654908.9915752509
$ws.Range("L3").Value = 936643.0641633343

$ws.Range("I4").Value = 911249.6041315054
$ws.Range("J4").Value = 1049684.73976003
$ws.Range("K4").Value = 935587.324069621
$ws.Range("L4").Value = 1071919.326292341

$ws.Range("I5").Value = 179929.469339055
$ws.Range("J5").Value = 357553.6942982785
$ws.Range("K5").Value = 206201.8971714219
$ws.Range("L5").Value = 381555.8017897776

$ws.Range("I6").Value = 246934.8822290114
$ws.Range("J6").Value = 428211.9206647582
$ws.Range("K6").Value = 274795.0882403546
$ws.Range("L6").Value = 453664.599179522

$ws.Range("I7").Value = -28459.16983042323
$ws.Range("J7").Value = 316261.1672119474
$ws.Range("K7").Value = -39625.69517414671
$ws.Range("L7").Value = 306059.5922636129

$ws.Range("I8").Value = 784943.8744965211
$ws.Range("J8").Value = 1053879.66278772
$ws.Range("K8").Value = 804303.5489011975
$ws.Range("L8").Value = 1071566.379410422

$ws.Range("I9").Value = 462666.2129662216
$ws.Range("J9").Value = 838746.6093900969
$ws.Range("K9").Value = 335119.0044015516
$ws.Range("L9").Value = 722221.336973579

$ws.Range("I10").Value = 4161.713205915598
$ws.Range("J10").Value = 258774.6683200491
$ws.Range("K10").Value = 7988.612946434011
$ws.Range("L10").Value = 262270.8681961013

$ws.Range("I11").Value = 333039.6142924197
$ws.Range("J11").Value = 618741.7211936882
$ws.Range("K11").Value = 355811.0408775005
$ws.Range("L11").Value = 639545.3649268548

$ws.Range("I12").Value = 5339047.169138927
$ws.Range("J12").Value = 5605829.629295891
$ws.Range("K12").Value = 5472499.694143863
$ws.Range("L12").Value = 5727749.91281968

$ws.Range("I13").Value = 6584437.838945085
$ws.Range("J13").Value = 7011485.666034247
$ws.Range("K13").Value = 7003038.973700793
$ws.Range("L13").Value = 7393913.567124221

$ws.Range("I14").Value = 756767.4396967822
$ws.Range("J14").Value = 951905.5567331358
$ws.Range("K14").Value = 807394.0444712352
$ws.Range("L14").Value = 998157.2856361279

$ws.Range("I15").Value = 680879.6059397161
$ws.Range("J15").Value = 886172.1217645985
$ws.Range("K15").Value = 739993.5414474566
$ws.Range("L15").Value = 940177.7524339625

$ws.Range("I16").Value = 267187.0574480658
$ws.Range("J16").Value = 610191.4122149925
$ws.Range("K16").Value = 283123.2977811521
$ws.Range("L16").Value = 624750.5293215535

$ws.Range("I17").Value = 4862420.20075058
$ws.Range("J17").Value = 5104014.293088192
$ws.Range("K17").Value = 5060732.987264478
$ws.Range("L17").Value = 5285189.967038705

$ws.Range("I18").Value = 18244.77936770647
$ws.Range("J18").Value = 351556.5911197374
$ws.Range("K18").Value = 21775.81392001805
$ws.Range("L18").Value = 354782.4928681996

$ws.Range("I19").Value = 141409.6994605679
$ws.Range("J19").Value = 339741.1544963364
$ws.Range("K19").Value = 148193.4257592498
$ws.Range("L19").Value = 345938.6680589496

$ws.Range("I20").Value = -16494.53329866743
$ws.Range("J20").Value = 186838.0016695139
$ws.Range("K20").Value = -30604.44252629921
$ws.Range("L20").Value = 173947.3940660673

$ws.Range("I21").Value = 236423.1070006753
$ws.Range("J21").Value = 418145.2731253565
$ws.Range("K21").Value = 269678.4096056995
$ws.Range("L21").Value = 448526.833325315

$ws.Range("I22").Value = 1165508.357386678
$ws.Range("J22").Value = 1634801.384678124
$ws.Range("K22").Value = 1273748.127808059
$ws.Range("L22").Value = 1733687.662759685

$ws.Range("I23").Value = 305560.2192948392
$ws.Range("J23").Value = 650982.9138272291
$ws.Range("K23").Value = 322655.30214271
$ws.Range("L23").Value = 666600.7325803621

$ws.Range("I24").Value = 2852475.818426631
$ws.Range("J24").Value = 3101908.486522131
$ws.Range("K24").Value = 3133639.714277288
$ws.Range("L24").Value = 3358775.727494132

$ws.Range("I25").Value = 747769.7197323454
$ws.Range("J25").Value = 981469.8207064534
$ws.Range("K25").Value = 772426.5731933251
$ws.Range("L25").Value = 1003995.962978066

$ws.Range("I26").Value = 1426825.30690537
$ws.Range("J26").Value = 1624356.430168747
$ws.Range("K26").Value = 1556702.134190955
$ws.Range("L26").Value = 1743010.008325046

$ws.Range("I27").Value = 596438.2140163889
$ws.Range("J27").Value = 815335.9143480232
$ws.Range("K27").Value = 654932.4084252042
$ws.Range("L27").Value = 868775.3585740097

$ws.Range("I28").Value = 422831.5099067052
$ws.Range("J28").Value = 447031.634790507
$ws.Range("K28").Value = 457440.5683490587
$ws.Range("L28").Value = 478649.9666102086

$ws.Range("I29").Value = 84231.64972324873
$ws.Range("J29").Value = 281713.8921881712
$ws.Range("K29").Value = 96141.21461400544
$ws.Range("L29").Value = 292594.2972466165

$ws.Range("I30").Value = 3787042.244251377
$ws.Range("J30").Value = 4042368.871669999
$ws.Range("K30").Value = 3888415.504299033
$ws.Range("L30").Value = 4134982.005863479

$ws.Range("I31").Value = 650385.7980818446
$ws.Range("J31").Value = 1073767.252441964
$ws.Range("K31").Value = 691566.270745585
$ws.Range("L31").Value = 1111389.132602127

$ws.Range("I32").Value = 602816.4688153502
$ws.Range("J32").Value = 851137.3783889906
$ws.Range("K32").Value = 640826.0584230279
$ws.Range("L32").Value = 885862.3859631785

$ws.Range("I33").Value = 5940850.042957444
$ws.Range("J33").Value = 6228782.625342092
$ws.Range("K33").Value = 5685373.795883409
$ws.Range("L33").Value = 5995383.246231748

$ws.Range("I34").Value = 3177123.64608264
$ws.Range("J34").Value = 3619030.517815856
$ws.Range("K34").Value = 3545962.381484392
$ws.Range("L34").Value = 3955996.215490545
